# Logged Week 16 and performed season sim from Week 17
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# YDS sheet - append this week's per-play yardage figures to the running logs
# ---------------------------------------------------------------------------
$wsYDS = $wb.Worksheets.Item("YDS")

$wsYDS.Range("B2").Value = $wsYDS.Range("B2").Value() + " 0 4 8 2 13 8 3 2 3 3 27 6 6 6 3 10 7 1 -1 6 5 4"
$wsYDS.Range("C2").Value = $wsYDS.Range("C2").Value() + " 1 1 1 27 11 17 3 8 3 12 5 8 13 8 15 13 3 3 12 8 30 4 10 4"
$wsYDS.Range("B3").Value = $wsYDS.Range("B3").Value() + " 1 8 4 33 4 6 11 34 1 9 2 -2 10 9 12 1 10 1 17 10 1 9 7 4"
$wsYDS.Range("C3").Value = $wsYDS.Range("C3").Value() + " 7 40 12 8 8 11 24 14 1 16 5 3 10 10 12 3 11 10 5 7"

# ---------------------------------------------------------------------------
# OFF sheet - down/distance tendency counters, updated with Week 16 totals
# ---------------------------------------------------------------------------
$wsOFF = $wb.Worksheets.Item("OFF")

$wsOFF.Range("C2").Value = 194
$wsOFF.Range("F2").Value = 50
$wsOFF.Range("G2").Value = 61
$wsOFF.Range("J2").Value = 20
$wsOFF.Range("L2").Value = 251
$wsOFF.Range("M2").Value = 173
$wsOFF.Range("O2").Value = 18
$wsOFF.Range("P2").Value = 11
$wsOFF.Range("Q2").Value = 449

$wsOFF.Range("B3").Value = 12
$wsOFF.Range("C3").Value = 153
$wsOFF.Range("D3").Value = 9
$wsOFF.Range("E3").Value = 22
$wsOFF.Range("F3").Value = 101
$wsOFF.Range("G3").Value = 43
$wsOFF.Range("H3").Value = 27
$wsOFF.Range("I3").Value = 55
$wsOFF.Range("J3").Value = 53

# ---------------------------------------------------------------------------
# DEF sheet - same tendency counters, defensive side
# ---------------------------------------------------------------------------
$wsDEF = $wb.Worksheets.Item("DEF")

$wsDEF.Range("B2").Value = 2
$wsDEF.Range("C2").Value = 168
$wsDEF.Range("D2").Value = 11
$wsDEF.Range("F2").Value = 48
$wsDEF.Range("G2").Value = 47
$wsDEF.Range("I2").Value = 5
$wsDEF.Range("J2").Value = 34
$wsDEF.Range("L2").Value = 260
$wsDEF.Range("M2").Value = 157
$wsDEF.Range("O2").Value = 26
$wsDEF.Range("P2").Value = 12
$wsDEF.Range("Q2").Value = 432

$wsDEF.Range("B3").Value = 11
$wsDEF.Range("C3").Value = 165
$wsDEF.Range("D3").Value = 8
$wsDEF.Range("E3").Value = 28
$wsDEF.Range("F3").Value = 104
$wsDEF.Range("G3").Value = 40
$wsDEF.Range("H3").Value = 14
$wsDEF.Range("I3").Value = 50
$wsDEF.Range("J3").Value = 56
$wsDEF.Range("N3").Value = 19

# ---------------------------------------------------------------------------
# ST sheet - special teams counters + per-game logs (kick/punt/FG distances)
# ---------------------------------------------------------------------------
$wsST = $wb.Worksheets.Item("ST")

$wsST.Range("B2").Value = 78
$wsST.Range("D2").Value = 48
$wsST.Range("F2").Value = 711
$wsST.Range("G2").Value = 694
$wsST.Range("J2").Value = 255
$wsST.Range("K2").Value = 232
$wsST.Range("B3").Value = 33

$wsST.Range("B4").Value = $wsST.Range("B4").Value() + " 64 63 52 65"
$wsST.Range("B5").Value = $wsST.Range("B5").Value() + " 24 28 16 15"
$wsST.Range("B6").Value = $wsST.Range("B6").Value() + " 21 27"
$wsST.Range("D3").Value = $wsST.Range("D3").Value() + " 45 44 56 49 41"
$wsST.Range("D4").Value = $wsST.Range("D4").Value() + " 4 0 15 0 8"
$wsST.Range("D5").Value = $wsST.Range("D5").Value() + " 9"

# ---------------------------------------------------------------------------
# TURNS sheet - turnovers
# ---------------------------------------------------------------------------
$wsTURNS = $wb.Worksheets.Item("TURNS")

$wsTURNS.Range("C2").Value = 12
